$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.222.58'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.55%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.005.73'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.18%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.66'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.90%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.625'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.56%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.72'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.23%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.388'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0813'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.17%  '
$ws.Range("E11").Value = '  +0.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.03'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.06%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.29'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.297.79'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.845'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.53%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.47'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.13%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.015.54'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.73%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.086.16'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.34%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.25'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.31%  '
$ws.Range("E20").Value = '  +1.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.21'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '230.54'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.49'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.45'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.32%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.32'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.12%  '
$ws.Range("E28").Value = '  -3.55%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.67'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.99%  '
$ws.Range("E30").Value = '  +11.93%  '
$ws.Range("E31").Value = '  +1.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.81'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.82%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0653'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.47'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.19%  '
$ws.Range("E35").Value = '  +6.45%  '
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.48'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.71%  '
$ws.Range("B37").Value = 'BinanceUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.08%  '
$ws.Range("E38").Value = '  +2.36%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.32'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.74%  '
$ws.Range("E40").Value = '  -0.20%  '
$ws.Range("E41").Value = '  +0.88%  '
$ws.Range("E42").Value = '  +0.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0214'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.93%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.56'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.93%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.12'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.80%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.371.25'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.29%  '
$ws.Range("E47").Value = '  +0.92%  '
$ws.Range("E48").Value = '  +2.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.07'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +12.74%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '46.77'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.36%  '
$ws.Range("E51").Value = '  -0.01%  '
